# weapons.xlsx update
# - fix bomb: base damage (基础伤害) for the row-6 weapon (E6) 20 -> 5
# - leaves the cell selection positioned at G6 (matches the saved session state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the bomb's base damage value
$ws.Range("E6").Value = 5

# Update the active selection to G6, as captured in the saved file
$ws.Range("G6").Select()
